# Apply cryptocurrency price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    # Force the cell to Text format so numeric-looking strings (e.g. "42.00")
    # are not reinterpreted/truncated as numbers by Excel, then restore the
    # default "Normal" style so no stray formatting is left behind.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "61.837.30"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "3.417.50"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  +0.01%  "
Set-TextValue $ws.Range("D5") "409.93"
$ws.Range("E5").Value = "  -0.01%  "
Set-TextValue $ws.Range("D6") "128.96"
$ws.Range("E6").Value = "  -0.68%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  +0.02%  "
Set-TextValue $ws.Range("D9") "0.735"
$ws.Range("E9").Value = "  -3.32%  "
$ws.Range("E10").Value = "  -0.50%  "
Set-TextValue $ws.Range("D11") "43.35"
$ws.Range("E11").Value = "  +0.18%  "
Set-TextValue $ws.Range("D12") "0.0000223"
$ws.Range("E12").Value = "  +14.50%  "
Set-TextValue $ws.Range("D13") "9.28"
$ws.Range("E13").Value = "  +5.56%  "
$ws.Range("D14").Value = "3.961.38"
$ws.Range("E14").Value = "  -0.47%  "
$ws.Range("E15").Value = "  +0.26%  "
Set-TextValue $ws.Range("D16") "21.24"
$ws.Range("E16").Value = "  +4.01%  "
$ws.Range("D17").Value = "3.422.29"
$ws.Range("E17").Value = "  +0.32%  "
Set-TextValue $ws.Range("D18") "12.34"
$ws.Range("E18").Value = "  +8.25%  "
$ws.Range("E19").Value = "  +3.00%  "
$ws.Range("D20").Value = "61.808.69"
$ws.Range("E20").Value = "  -0.52%  "
Set-TextValue $ws.Range("D21") "492.95"
$ws.Range("E21").Value = "  +30.23%  "
Set-TextValue $ws.Range("D22") "91.61"
$ws.Range("E22").Value = "  +4.17%  "
Set-TextValue $ws.Range("D23") "3.34"
$ws.Range("E23").Value = "  +4.90%  "
Set-TextValue $ws.Range("D24") "13.49"
$ws.Range("E24").Value = "  +0.44%  "
Set-TextValue $ws.Range("D25") "3.33"
$ws.Range("E25").Value = "  +3.83%  "
Set-TextValue $ws.Range("D26") "34.39"
$ws.Range("E26").Value = "  +8.28%  "
Set-TextValue $ws.Range("D27") "9.29"
$ws.Range("E27").Value = "  +9.99%  "
$ws.Range("E28").Value = "  -1.18%  "
Set-TextValue $ws.Range("D29") "12.14"
$ws.Range("E29").Value = "  +2.56%  "
$ws.Range("E30").Value = "  -1.75%  "
Set-TextValue $ws.Range("D31") "0.115"
$ws.Range("E31").Value = "  -1.52%  "
$ws.Range("E32").Value = "  -2.16%  "
Set-TextValue $ws.Range("D33") "42.00"
$ws.Range("E33").Value = "  -4.61%  "
$ws.Range("E34").Value = "  +0.07%  "
Set-TextValue $ws.Range("D35") "58.54"
$ws.Range("E35").Value = "  +12.20%  "
$ws.Range("E36").Value = "  +1.05%  "
Set-TextValue $ws.Range("D37") "0.998"
$ws.Range("E37").Value = "  +0.09%  "
Set-TextValue $ws.Range("D38") "3.45"
$ws.Range("E38").Value = "  +2.61%  "
$ws.Range("E39").Value = "  +3.55%  "
Set-TextValue $ws.Range("D40") "2.72"
$ws.Range("E40").Value = "  +17.28%  "
Set-TextValue $ws.Range("D41") "147.92"
$ws.Range("E41").Value = "  +3.65%  "
$ws.Range("E42").Value = "  +0.22%  "
Set-TextValue $ws.Range("D43") "0.319"
$ws.Range("E43").Value = "  +1.39%  "
Set-TextValue $ws.Range("D44") "2.09"
$ws.Range("E44").Value = "  +5.88%  "
Set-TextValue $ws.Range("D45") "4.35"
$ws.Range("E45").Value = "  +8.63%  "
$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue $ws.Range("D46") "16.79"
$ws.Range("E46").Value = "  +0.47%  "
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue $ws.Range("D47") "2.33"
$ws.Range("E47").Value = "  +20.01%  "
Set-TextValue $ws.Range("D48") "23.01"
$ws.Range("E48").Value = "  +5.55%  "
Set-TextValue $ws.Range("D49") "118.32"
$ws.Range("E49").Value = "  +26.62%  "
$ws.Range("E50").Value = "  +17.33%  "
$ws.Range("D51").Value = "2.141.20"
$ws.Range("E51").Value = "  +1.13%  "
